# ----------------------------------------------------------------------
# edit.ps1
# Reproduces the "Add files via upload" commit:
#  - duplicates "resumen" into a new "conv" sheet (freezing its old
#    array-formula values into plain numbers) and adds a summary row/cols
#  - adds row 30 to "bets" (second betting "base")
#  - adds row 4 to "resumen" (now recalculated against the new row 30)
#  - updates sheet views / active tab
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$bets    = $wb.Worksheets.Item("bets")
$resumen = $wb.Worksheets.Item("resumen")

# ------------------------------------------------------------------
# 1) Duplicate "resumen" -> "conv" (placed after "resumen", i.e. last)
#    This captures resumen's CURRENT (pre row-30) array-formula values.
# ------------------------------------------------------------------
$resumen.Copy($null, $resumen)
$conv = $wb.Worksheets.Item($wb.Worksheets.Count)
$conv.Name = "conv"

# Freeze the old array formulas in B3:F3 into plain literal values
# (still based on the old M29 = 2.575 baseline).
$conv.Range("B3").Value = 474347.83
$conv.Range("C3").Value = 56162.951294999999
$conv.Range("D3").Value = 83730.495419999992
$conv.Range("E3").Value = 33743.913115000003
$conv.Range("F3").Value = 26457.826229999999

# Column G: running totals per row
$conv.Range("G2").Formula = "=SUM(B2:F2)"
$conv.Range("G3").Formula = "=SUM(B3:F3)"

# New row 4: a fresh "base" (+1,000,000 on top of the accumulated bank)
$conv.Range("A4").Value = 3
$conv.Range("B4").Formula = "=B3+1000000"
$conv.Range("C4").Value = 56162.951294999999
$conv.Range("D4").Value = 83730.495419999992
$conv.Range("E4").Value = 33743.913115000003
$conv.Range("F4").Value = 26457.826229999999
$conv.Range("G4").Formula = "=SUM(B4:F4)"

# Extra helper columns on row 4 (percentage -> unit conversion table)
$conv.Range("H4").Formula = "=2.575+100"
$conv.Range("I4").Formula = "=G4/100"
$conv.Range("J4").Formula = "=I4*2"
$conv.Range("K4").Formula = "=I4+J4"
$conv.Range("L4").Formula = "=K4+J4"

# Header row additions (G1, I1:L1)
$conv.Range("G1").Value  = "BASE"
$conv.Range("I1").Value  = "1U"
$conv.Range("J1").Value  = "2U"
$conv.Range("K1").Value  = "3U"
$conv.Range("L1").Value  = "5U"

$conv.Range("A4:F4").Select()
$conv.Range("A1").Select()

# ------------------------------------------------------------------
# 2) "bets": add row 30 (a new bet entry) - copy formatting from row 29
# ------------------------------------------------------------------
$bets.Range("A29:M29").Copy()
$bets.Range("A30:M30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$bets.Range("A30").Value = 29
$bets.Range("B30").Value = 45224
$bets.Range("C30").Value = 1
$bets.Range("D30").Value = 1674443.0160600001
$bets.Range("E30").Value = 13
$bets.Range("F30").Formula = "=D30+E30"
$bets.Range("G30").Value = "ESPORTS"
$bets.Range("H30").Value = "IBERIAN CUP"
$bets.Range("I30").Value = "LH"
$bets.Range("J30").Value = "GANA SERIE"
$bets.Range("K30").Value = 1
$bets.Range("L30").Value = 0
$bets.Range("M30").Formula = "=ROUND((F30/`$D`$30-1)*100, 3)+`$M`$29"

# ------------------------------------------------------------------
# 3) "resumen": add row 4 (literal snapshot of conv!row4) + widen cols
# ------------------------------------------------------------------
$resumen.Range("A4").Value = 3
$resumen.Range("B4").Value = 1474347.83
$resumen.Range("C4").Value = 56162.951294999999
$resumen.Range("D4").Value = 83730.495419999992
$resumen.Range("E4").Value = 33743.913115000003
$resumen.Range("F4").Value = 26457.826229999999

$resumen.Range("B2:F2").Select()
$resumen.Range("H6").Select()

# ------------------------------------------------------------------
# 4) Activate "resumen" (becomes the active tab) and select cells to
#    match the final view state recorded in the workbook.
# ------------------------------------------------------------------
$bets.Activate()
$bets.Range("H21").Select()

$resumen.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 2
$resumen.Range("H6").Select()

$conv.Activate()
$conv.Range("A4:F4").Select()

$resumen.Activate()
